# Remove the row for account 005000460 / MARIANA / 3000 from the "Export" sheet.
# This is row 6 (row 1 is the header: Conta, Nome, Saldo).
# Deleting the entire row shifts all subsequent rows up by one, matching the
# target diff which removes that <x:row> block entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Export")

$ws.Rows.Item(6).Delete()
